$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (people interested) counts in column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 31
$wsExpo.Range("F4").Value = 60
$wsExpo.Range("F5").Value = 5049
$wsExpo.Range("F8").Value = 305
$wsExpo.Range("F9").Value = 53

# Sheet "全部类型" (All types) - same events appear here with the same updated counts
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 31
$wsAll.Range("F8").Value = 60
$wsAll.Range("F9").Value = 5049
$wsAll.Range("F13").Value = 305
$wsAll.Range("F14").Value = 53
